$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.699.41"
$ws.Range("E2").Value = "  -6.99%  "
$ws.Range("D3").Value = "2.591.10"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "301.48"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "96.57"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.561"
$ws.Range("E9").Value = "  -4.52%  "
$ws.Range("D10").Value = "36.93"
$ws.Range("E10").Value = "  -8.45%  "
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("D12").Value = "7.83"
$ws.Range("E12").Value = "  -5.78%  "
$ws.Range("D13").Value = "2.984.78"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D15").Value = "2.581.50"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "0.894"
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("D17").Value = "14.37"
$ws.Range("E17").Value = "  -5.31%  "
$ws.Range("D18").Value = "43.668.88"
$ws.Range("E18").Value = "  -7.22%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0983"
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("D21").Value = "12.39"
$ws.Range("E21").Value = "  -6.31%  "
$ws.Range("D22").Value = "73.32"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "265.97"
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "29.28"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "10.31"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "38.16"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  -5.22%  "
$ws.Range("D32").Value = "3.61"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").Value = "2.22"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").Value = "151.96"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0820"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  -5.49%  "
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").Value = "24.34"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").Value = "16.98"
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").Value = "0.0317"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("D43").Value = "3.89"
$ws.Range("E43").Value = "  -6.18%  "
$ws.Range("D44").Value = "2.033.90"
$ws.Range("E44").Value = "  -4.43%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "88.24"
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("D47").Value = "9.14"
$ws.Range("E47").Value = "  -4.87%  "
$ws.Range("D48").Value = "1.61"
$ws.Range("E48").Value = "  +4.37%  "
$ws.Range("D49").Value = "2.838.01"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").Value = "105.93"
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("D51").Value = "0.192"
$ws.Range("E51").Value = "  -5.56%  "
